# Add 2022-Q1 data: insert a new "2022-Q1" sheet (cloning the layout/style of the
# most recent quarterly sheet) right before the "总计" (totals) summary sheet, and
# update the "总计" sheet with a new leading row for 2022-Q1, shifting the older
# rows down.
#
# NOTE: worksheet references returned/stored before a sheet-order-changing call
# (like Worksheets.Add) can go stale - they resolve by position, so after
# inserting a sheet we always re-fetch the sheets we need by name.

function Set-TextValue($range, [string]$val) {
    # Force a value to be stored as text, even if it looks like a pure number
    # (e.g. "070031", "0.60", "0.0187"), matching the original inlineStr cells.
    $range.NumberFormat = "@"
    $range.Value = $val
    # Drop back to the default "Normal" style so we don't leave a stray
    # number-format-only style behind (keeps cellXfs clean, like the source).
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# --- locate the existing sheets -------------------------------------------------
$q4sheet = $wb.Worksheets.Item("2021-Q4")  # used as a style/layout template
$totalSheetBefore = $wb.Worksheets.Item("总计")

# --- 1. insert the new "2022-Q1" sheet right before "总计" ----------------------
$q1sheet = $wb.Worksheets.Add($totalSheetBefore)
$q1sheet.Name = "2022-Q1"

# Re-fetch by name: inserting a sheet shifts positions, so any reference
# captured before this point (e.g. $totalSheetBefore) can now point at the
# wrong slot.
$q1sheet = $wb.Worksheets.Item("2022-Q1")
$totalSheet = $wb.Worksheets.Item("总计")

# Clone header row + first data-row formatting from the "2021-Q4" sheet so the
# new sheet keeps the same fund-table look (bold/centered/bordered header, etc.)
$q4sheet.Range("B1:H2").Copy($q1sheet.Range("B1"))
$q4sheet.Range("A2").Copy($q1sheet.Range("A2"))

# Overwrite with the 2022-Q1 fund data
Set-TextValue $q1sheet.Range("B2") "070031"
Set-TextValue $q1sheet.Range("C2") "嘉实全球房地产(QDII)"
Set-TextValue $q1sheet.Range("D2") "0.60"
Set-TextValue $q1sheet.Range("E2") "95.08"
Set-TextValue $q1sheet.Range("F2") "3.12"
Set-TextValue $q1sheet.Range("G2") "0.0187"
$q1sheet.Range("H2").Value = 6

# --- 2. update the "总计" sheet: add a 2022-Q1 row on top, push the rest down ---
# Clone row 2's formatting down into the new row 5 first, so every row keeps the
# same per-cell styling (style on column A) as the existing rows.
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A5:D5"))

# Row 5 (was conceptually row 4 / "2021-Q2")
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 1
$totalSheet.Range("D5").Value = 0.02

# Row 4 (was row 3 / "2021-Q3")
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.01

# Row 3 (was row 2 / "2021-Q4")
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.02

# Row 2 (new / "2022-Q1")
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.02

Write-Host "Done: added 2022-Q1 sheet and updated 总计 totals"
